$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp shown in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 10:22"

# --- Rows 30/31: Polonia overtakes Arabia Saudita (sorted by Casos totales desc) ---
# Row 30 becomes Polonia with refreshed stats
$ws.Range("A30").Value = "Polonia"
$ws.Range("B30").Value = 9453
$ws.Range("C30").Value = 166
$ws.Range("D30").Value = 1133
$ws.Range("E30").Value = 7958
$ws.Range("F30").Value = 160
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = 362

# Row 31 becomes Arabia Saudita, carrying the prior (unrefreshed) stats of the old row 30
$ws.Range("A31").Value = "Arabia Saudita"
$ws.Range("B31").Value = 9362
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = 1398
$ws.Range("E31").Value = 7867
$ws.Range("F31").Value = 97
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 97

# --- Row 36: Dinamarca refreshed stats ---
$ws.Range("B36").Value = 7515
$ws.Range("C36").Value = 131
$ws.Range("E36").Value = 3019

# --- Rows 42/43: Filipinas overtakes Serbia (sorted by Casos totales desc) ---
# Row 42 becomes Filipinas with refreshed stats
$ws.Range("A42").Value = "Filipinas"
$ws.Range("B42").Value = 6459
$ws.Range("C42").Value = 200
$ws.Range("D42").Value = 613
$ws.Range("E42").Value = 5418
$ws.Range("F42").Value = 1
$ws.Range("G42").Value = 19
$ws.Range("H42").Value = 428

# Row 43 becomes Serbia, carrying the prior (unrefreshed) stats of the old row 42
$ws.Range("A43").Value = "Serbia"
$ws.Range("B43").Value = 6318
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 753
$ws.Range("E43").Value = 5443
$ws.Range("F43").Value = 120
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 122

# --- Row 67: Kazajistan refreshed stats (D, E only) ---
$ws.Range("D67").Value = 417
$ws.Range("E67").Value = 1321

# --- Row 113: refreshed stats (D, E only) ---
$ws.Range("D113").Value = 27
$ws.Range("E113").Value = 275

# --- Row 116: refreshed stats (D, E only) ---
$ws.Range("D116").Value = 97
$ws.Range("E116").Value = 191
